$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting existing rows 22:50 down to 23:51.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly price record.
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "Vega Modelo de Temuco"
$ws.Range("C22").Value = "La Araucanía"
$ws.Range("D22").Value = 44868
$ws.Range("D22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 100112042
$ws.Range("G22").Value = "Locoto"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 2500
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = 2500
$ws.Range("N22").Value = "$/kilo"
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 2500
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = "Hortaliza"
